$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value that was updated from 2023-09-03
# (serial 45172) to 2023-09-06 (serial 45175) for every data row (rows 2-236).
$ws.Range("C2:C236").Value = 45175
